# Daily attendance processing
# "Recorded By" (column G) lists the accounts that touched a session's
# attendance record as a comma-separated string, e.g. "System, user@x.com".
# The daily sync re-orders that list so the System marker that kicked off
# the record is rotated to the back, keeping the human/automation accounts
# that follow up front: "a, b, c" -> "b, c, a".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $text = $cell.Text

    if ($text -and $text.Length -gt 0) {
        $parts = $text -split ", "

        if ($parts.Length -gt 1) {
            # rotate left by one: move the first entry to the end
            $rotated = $parts[1..($parts.Length - 1)] + $parts[0]
            $cell.Value = $rotated -join ", "
        }
    }
}
